$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.NumberFormat = "General"
    $rng.Style = "Normal"
}

# Row 2
Set-TextValue 'D2' '30.038.63'
Set-TextValue 'E2' '  +0.67%  '

# Row 3
Set-TextValue 'D3' '1.885.10'
Set-TextValue 'E3' '  -0.13%  '

# Row 4
Set-TextValue 'D4' '1.001'
Set-TextValue 'E4' '  +0.05%  '

# Row 5
Set-TextValue 'D5' '0.7365'
Set-TextValue 'E5' '  -2.09%  '

# Row 6
Set-TextValue 'D6' '242.33'
Set-TextValue 'E6' '  +0.06%  '

# Row 7
Set-TextValue 'D7' '1.001'
Set-TextValue 'E7' '  +0.07%  '

# Row 8
Set-TextValue 'D8' '0.3162'
Set-TextValue 'E8' '  +1.09%  '

# Row 9
Set-TextValue 'D9' '0.07177'
Set-TextValue 'E9' '  +0.81%  '

# Row 10
Set-TextValue 'D10' '24.64'
Set-TextValue 'E10' '  -2.58%  '

# Row 11
Set-TextValue 'E11' '  -2.00%  '

# Row 12
Set-TextValue 'D12' '0.7562'
Set-TextValue 'E12' '  -0.48%  '

# Row 13
Set-TextValue 'D13' '1.908.38'
Set-TextValue 'E13' '  +1.74%  '

# Row 14
Set-TextValue 'D14' '5.404'
Set-TextValue 'E14' '  +0.65%  '

# Row 15
Set-TextValue 'D15' '92.54'
Set-TextValue 'E15' '  -0.80%  '

# Row 16
Set-TextValue 'D16' '6.143'
Set-TextValue 'E16' '  +0.22%  '

# Row 17
Set-TextValue 'D17' '30.031.85'
Set-TextValue 'E17' '  +0.42%  '

# Row 18
Set-TextValue 'D18' '249.75'
Set-TextValue 'E18' '  +2.74%  '

# Row 19
Set-TextValue 'D19' '13.55'
Set-TextValue 'E19' '  -1.16%  '

# Row 20
Set-TextValue 'D20' '0.000007844'
Set-TextValue 'E20' '  +0.07%  '

# Row 21
Set-TextValue 'D21' '2.153.74'
Set-TextValue 'E21' '  +0.52%  '

# Row 22
Set-TextValue 'D22' '0.9994'
Set-TextValue 'E22' '  -0.05%  '

# Row 23
Set-TextValue 'D23' '1.001'
Set-TextValue 'E23' '  +0.03%  '

# Row 24
Set-TextValue 'D24' '7.885'
Set-TextValue 'E24' '  -1.36%  '

# Row 25
Set-TextValue 'D25' '0.1568'
Set-TextValue 'E25' '  -1.44%  '

# Row 26
Set-TextValue 'D26' '9.262'
Set-TextValue 'E26' '  -1.11%  '

# Row 27
Set-TextValue 'D27' '164.14'
Set-TextValue 'E27' '  +0.65%  '

# Row 28
Set-TextValue 'D28' '18.64'
Set-TextValue 'E28' '  -0.43%  '

# Row 29
Set-TextValue 'D29' '2.045'
Set-TextValue 'E29' '  +0.73%  '

# Row 30
Set-TextValue 'D30' '1.473'
Set-TextValue 'E30' '  -0.61%  '

# Row 31
Set-TextValue 'D31' '4.558'
Set-TextValue 'E31' '  +1.01%  '

# Row 32
Set-TextValue 'D32' '1.532'
Set-TextValue 'E32' '  +0.05%  '

# Row 33
Set-TextValue 'D33' '4.184'
Set-TextValue 'E33' '  +0.60%  '

# Row 34
Set-TextValue 'D34' '0.05329'
Set-TextValue 'E34' '  -1.76%  '

# Row 35
Set-TextValue 'D35' '1.247'
Set-TextValue 'E35' '  +0.54%  '

# Row 36
Set-TextValue 'D36' '0.7666'
Set-TextValue 'E36' '  +1.90%  '

# Row 37
Set-TextValue 'D37' '0.9996'
Set-TextValue 'E37' '  -0.42%  '

# Row 38
Set-TextValue 'D38' '2.731'
Set-TextValue 'E38' '  +0.75%  '

# Row 39
Set-TextValue 'D39' '0.01955'
Set-TextValue 'E39' '  +0.44%  '

# Row 40
Set-TextValue 'D40' '2.759'
Set-TextValue 'E40' '  -0.48%  '

# Row 41
Set-TextValue 'D41' '0.4553'
Set-TextValue 'E41' '  +1.87%  '

# Row 42
Set-TextValue 'D42' '1.093.45'
Set-TextValue 'E42' '  -0.75%  '

# Row 43
Set-TextValue 'D43' '6.048'
Set-TextValue 'E43' '  -1.02%  '

# Row 44
Set-TextValue 'B44' 'TrustWalletToken'
Set-TextValue 'C44' 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue 'D44' '0.8810'
Set-TextValue 'E44' '  +2.42%  '

# Row 45
Set-TextValue 'B45' 'Aave'
Set-TextValue 'C45' 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue 'D45' '72.27'
Set-TextValue 'E45' '  -0.25%  '

# Row 46
Set-TextValue 'B46' 'Quant'
Set-TextValue 'C46' 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextValue 'D46' '104.26'
Set-TextValue 'E46' '  +1.69%  '

# Row 47
Set-TextValue 'B47' 'PaxDollar'
Set-TextValue 'C47' 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
Set-TextValue 'D47' '1.002'
Set-TextValue 'E47' '  +0.19%  '

# Row 48
Set-TextValue 'D48' '1.851'
Set-TextValue 'E48' '  -0.44%  '

# Row 49
Set-TextValue 'D49' '7.532'
Set-TextValue 'E49' '  -2.59%  '

# Row 50
Set-TextValue 'B50' 'EnergySwap'
Set-TextValue 'C50' 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue 'D50' '9.532'
Set-TextValue 'E50' '  -2.33%  '

# Row 51
Set-TextValue 'B51' 'RocketPoolETH'
Set-TextValue 'C51' 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
Set-TextValue 'D51' '2.053.20'
Set-TextValue 'E51' '  +0.56%  '
